$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 17
$ws.Range("F4").Value = 575
$ws.Range("F7").Value = 1918
$ws.Range("F8").Value = 5203
$ws.Range("F9").Value = 1453
$ws.Range("F11").Value = 3053
$ws.Range("F14").Value = 1268
$ws.Range("F15").Value = 4164
$ws.Range("F16").Value = 985
$ws.Range("F18").Value = 1637
$ws.Range("F19").Value = 2580
$ws.Range("F20").Value = 25
$ws.Range("F21").Value = 14
$ws.Range("F22").Value = 119
$ws.Range("F24").Value = 951
$ws.Range("F26").Value = 75
$ws.Range("E27").Value = "2024.05.19 10:00-05.19 17:00"
$ws.Range("F27").Value = 72
$ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202404/3WNfBWY61713863269103.jpeg"
$ws.Range("F28").Value = 199
$ws.Range("F30").Value = 343
$ws.Range("F31").Value = 22
$ws.Range("F32").Value = 113
$ws.Range("F34").Value = 212
$ws.Range("F35").Value = 1614
$ws.Range("F36").Value = 2141
$ws.Range("F37").Value = 997
$ws.Range("F39").Value = 242
$ws.Range("F40").Value = 590
$ws.Range("F41").Value = 257
$ws.Range("F42").Value = 84
$ws.Range("F44").Value = 386
$ws.Range("F45").Value = 289
$ws.Range("F47").Value = 126

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 27
$ws.Range("G2").Value = 578
$ws.Range("F6").Value = 11
$ws.Range("F7").Value = 32
$ws.Range("F8").Value = 11

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 720

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 720
$ws.Range("F3").Value = 27
$ws.Range("G3").Value = 578
$ws.Range("F4").Value = 17
$ws.Range("F6").Value = 575
$ws.Range("F8").Value = 1918
$ws.Range("F9").Value = 5203
$ws.Range("F10").Value = 1453
$ws.Range("F13").Value = 3053
$ws.Range("F15").Value = 1268
$ws.Range("F16").Value = 4164
$ws.Range("F17").Value = 985
$ws.Range("F18").Value = 1637
$ws.Range("F19").Value = 11
$ws.Range("F20").Value = 2580
$ws.Range("F21").Value = 32
$ws.Range("F22").Value = 25
$ws.Range("F23").Value = 11
$ws.Range("F27").Value = 951
$ws.Range("E29").Value = "2024.05.19 10:00-05.19 17:00"
$ws.Range("F29").Value = 72
$ws.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202404/3WNfBWY61713863269103.jpeg"
$ws.Range("F30").Value = 199
$ws.Range("F33").Value = 343
$ws.Range("F34").Value = 22
$ws.Range("F36").Value = 1614
$ws.Range("F37").Value = 2141
$ws.Range("F38").Value = 997
$ws.Range("F41").Value = 242
$ws.Range("F42").Value = 590
$ws.Range("F43").Value = 257
$ws.Range("F45").Value = 386
$ws.Range("F46").Value = 289
$ws.Range("F48").Value = 126
